$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text-ish "Price" value into column D without letting Excel
# auto-coerce numeric-looking strings ("0.9994", "240.18", ...) into numbers -
# the source data keeps these as plain text even though they look numeric.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# --- Price / Volume(1h) refresh for rows 2-45 (Coin/Link unchanged) ---
$updates = @(
    @{ Row = 2; D = '29.332.06'; E = '  -0.04%  ' }
    @{ Row = 3; D = '1.843.36'; E = '  -0.20%  ' }
    @{ Row = 4; D = '0.9994'; E = '  +0.15%  ' }
    @{ Row = 5; D = '240.18'; E = '  -0.06%  ' }
    @{ Row = 6; D = '0.6270'; E = '  -0.63%  ' }
    @{ Row = 7; D = $null; E = '  +0.24%  ' }
    @{ Row = 8; D = '0.07399'; E = '  -2.40%  ' }
    @{ Row = 9; D = $null; E = '  -0.98%  ' }
    @{ Row = 10; D = '24.76'; E = '  +1.46%  ' }
    @{ Row = 11; D = '0.07741'; E = '  +0.03%  ' }
    @{ Row = 12; D = '1.842.53'; E = '  -0.24%  ' }
    @{ Row = 13; D = '4.984'; E = '  -0.53%  ' }
    @{ Row = 14; D = '0.6773'; E = '  -0.29%  ' }
    @{ Row = 15; D = '0.00001016'; E = '  -2.87%  ' }
    @{ Row = 16; D = '82.02'; E = '  -1.41%  ' }
    @{ Row = 17; D = '6.260'; E = '  +2.19%  ' }
    @{ Row = 18; D = '29.363.95'; E = '  +0.06%  ' }
    @{ Row = 19; D = '228.78'; E = '  -0.35%  ' }
    @{ Row = 20; D = $null; E = '  -0.33%  ' }
    @{ Row = 21; D = '1.000'; E = '  +0.16%  ' }
    @{ Row = 22; D = '7.441'; E = '  -0.13%  ' }
    @{ Row = 23; D = $null; E = '  +0.17%  ' }
    @{ Row = 24; D = '158.73'; E = '  +0.15%  ' }
    @{ Row = 25; D = '8.464'; E = '  +0.26%  ' }
    @{ Row = 26; D = '0.1351'; E = '  -3.14%  ' }
    @{ Row = 27; D = '17.42'; E = '  -1.23%  ' }
    @{ Row = 28; D = '0.06620'; E = '  +16.51%  ' }
    @{ Row = 29; D = '1.459'; E = '  +2.54%  ' }
    @{ Row = 30; D = $null; E = '  +1.01%  ' }
    @{ Row = 31; D = '4.069'; E = '  -1.20%  ' }
    @{ Row = 32; D = '4.068'; E = '  +0.64%  ' }
    @{ Row = 33; D = '1.835'; E = '  +0.55%  ' }
    @{ Row = 34; D = '1.138'; E = '  -1.43%  ' }
    @{ Row = 35; D = $null; E = '  -0.83%  ' }
    @{ Row = 36; D = '2.575'; E = '  -0.09%  ' }
    @{ Row = 37; D = '0.01855'; E = '  +1.64%  ' }
    @{ Row = 38; D = '2.819'; E = '  +3.65%  ' }
    @{ Row = 39; D = '1.243.46'; E = '  +0.02%  ' }
    @{ Row = 40; D = $null; E = '  +5.79%  ' }
    @{ Row = 41; D = '0.9365'; E = '  +3.79%  ' }
    @{ Row = 42; D = '1.001'; E = '  +0.17%  ' }
    @{ Row = 43; D = '2.031.42'; E = '  +1.23%  ' }
    @{ Row = 44; D = '100.60'; E = '  -0.92%  ' }
    @{ Row = 45; D = '65.59'; E = '  -0.18%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) { Set-TextValue $ws.Cells.Item($u.Row, 4) $u.D }
    if ($null -ne $u.E) { $ws.Cells.Item($u.Row, 5).Value = $u.E }
}

# --- Rows 46-51: two coin pairs swap places (re-ranked) + EnergySwap refresh ---

# Row 46: RenderToken -> Aptos
$ws.Cells.Item(46, 2).Value = 'Aptos'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Cells.Item(46, 4) '7.034'
$ws.Cells.Item(46, 5).Value = '  -1.38%  '

# Row 47: Aptos -> RenderToken
$ws.Cells.Item(47, 2).Value = 'RenderToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Cells.Item(47, 4) '1.709'
$ws.Cells.Item(47, 5).Value = '  +2.13%  '

# Row 48: EnergySwap (Coin/Link unchanged) - Price/Volume refresh
Set-TextValue $ws.Cells.Item(48, 4) '9.018'
$ws.Cells.Item(48, 5).Value = '  -0.04%  '

# Row 49: BabyDogeCoin -> Algorand
$ws.Cells.Item(49, 2).Value = 'Algorand'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Cells.Item(49, 4) '0.1150'
$ws.Cells.Item(49, 5).Value = '  -1.28%  '

# Row 50: Algorand -> TheSandbox
$ws.Cells.Item(50, 2).Value = 'TheSandbox'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Cells.Item(50, 4) '0.3906'
$ws.Cells.Item(50, 5).Value = '  -1.12%  '

# Row 51: TheSandbox -> BabyDogeCoin
$ws.Cells.Item(51, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Cells.Item(51, 4) '0.00000000113'
$ws.Cells.Item(51, 5).Value = '  -1.93%  '
